# Applies hybrid bold+color highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across the resume body.
#
# Strategy: for each target paragraph, locate its full (pre-edit) text via
# Find.Execute to get its absolute Start offset in the document, then for
# each metric substring (in left-to-right order) compute the absolute
# sub-range and apply Bold + the brand accent color (#2C3E50) to it. This
# naturally splits the paragraph's single run into the run sequence shown
# in the diff, since Word COM copies the existing run formatting onto each
# newly materialized run and only the explicitly-touched sub-ranges gain
# the Bold/Color overrides.

$accentColor = 5258796  # BGR-encoded #2C3E50 (R=0x2C,G=0x3E,B=0x50)

function Set-MetricHighlights($fullText, $segments) {
    $doc = $word.ActiveDocument
    $searchRange = $doc.Content
    $ok = $searchRange.Find.Execute($fullText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $fullText"
        return
    }
    $baseStart = $searchRange.Start
    $searchFrom = 0
    foreach ($seg in $segments) {
        $idx = $fullText.IndexOf($seg, $searchFrom)
        if ($idx -lt 0) {
            Write-Output "SEGMENT NOT FOUND: $seg"
            continue
        }
        $absStart = $baseStart + $idx
        $absEnd = $absStart + $seg.Length
        $segRange = $doc.Range($absStart, $absEnd)
        $segRange.Font.Bold = $true
        $segRange.Font.Color = $accentColor
        $searchFrom = $idx + $seg.Length
    }
}

# --- Professional Experience bullets ---------------------------------

Set-MetricHighlights '• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%' @('23%', '64%')

Set-MetricHighlights '• Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes' @('±4.2%', '±2.1%', '71%', '87%')

Set-MetricHighlights '• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis' @('73.5%', '$4.7M')

Set-MetricHighlights '• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion' @('$2')

Set-MetricHighlights '• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%' @('57%')

# --- Key Achievements and Impact bullets ------------------------------

Set-MetricHighlights '• 178% accuracy improvement in racial classification algorithms' @('178%')

Set-MetricHighlights '• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%' @('73.5%')

Set-MetricHighlights '• $4.7M savings enabled nonprofit access' @('$4.7M')

Set-MetricHighlights '• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations' @('12,847')

Set-MetricHighlights '• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%' @('±4.2%', '±2.1%')

Set-MetricHighlights '• Increased voter turnout prediction accuracy from 71% to 87%' @('71%', '87%')
